# Updates the "cryptos" symbol list sheet with refreshed Price (column D)
# and Hora (column G) values, as produced by the scheduled GitHub Actions
# scraper run on Thu Dec 15 18:11:50 UTC 2022.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Price (column D) values, keyed by row number. Rows whose price is
# "--" (no data) in the source sheet are intentionally left untouched.
$priceUpdates = @{
    2  = "261.46"
    3  = "22.71"
    4  = "6.202"
    5  = "0.06097"
    6  = "3.511"
    7  = "6.708"
    8  = "1.358"
    9  = "0.7990"
    11 = "0.08133"
    13 = "0.03151"
    14 = "0.09263"
    15 = "3.906"
    16 = "0.001711"
    17 = "0.04822"
    18 = "0.0006205"
    19 = "0.006214"
    20 = "0.001099"
    21 = "0.003370"
    23 = "3.696"
    24 = "2.297"
    25 = "0.3373"
    26 = "0.1251"
    27 = "0.0006170"
    40 = "0.04588"
    41 = "0.1119"
    42 = "0.003133"
    43 = "0.003424"
    44 = "0.01023"
    45 = "0.002972"
    46 = "0.00006038"
    47 = "0.00000000751"
    48 = "0.7006"
    49 = "0.1387"
    50 = "0.00002102"
    51 = "0.01011"
}

# Every data row (2-51) gets its "Hora" (column G) timestamp bumped from
# 17 to 18.
$newHora = "18"

for ($row = 2; $row -le 51; $row++) {
    if ($priceUpdates.ContainsKey($row)) {
        # Leading apostrophe forces these to stay text cells (matching the
        # inline-string storage already used throughout the sheet) instead
        # of being reinterpreted as numbers, which would lose trailing
        # zeros / exact formatting (e.g. "0.7990", "1.358").
        $ws.Range("D$row").Value = "'" + $priceUpdates[$row]
    }
    $ws.Range("G$row").Value = "'" + $newHora
}
